# Add new columns I (I0) and J (IF) to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, centered, bordered) from H1 onto the
# new header cells so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for the new I and J columns, keyed by row number (2-37)
$iValues = @{
    2 = 1;  3 = 1;  4 = 1;  5 = 1;  6 = 1;  7 = 1;  8 = 1;  9 = 1; 10 = 1;
   11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1; 19 = 1;
   20 = 1; 21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1; 26 = 1; 27 = 1; 28 = 1;
   29 = 1; 30 = 5; 31 = 4; 32 = 1; 33 = 1; 34 = 1; 35 = 1; 36 = 1; 37 = 3
}
$jValues = @{
    2 = 7;  3 = 5;  4 = 5;  5 = 6;  6 = 4;  7 = 7;  8 = 6;  9 = 4; 10 = 4;
   11 = 2; 12 = 6; 13 = 5; 14 = 5; 15 = 7; 16 = 6; 17 = 6; 18 = 6; 19 = 6;
   20 = 4; 21 = 5; 22 = 5; 23 = 5; 24 = 5; 25 = 5; 26 = 5; 27 = 6; 28 = 5;
   29 = 6; 30 = 7; 31 = 7; 32 = 3; 33 = 7; 34 = 5; 35 = 4; 36 = 3; 37 = 4
}

foreach ($row in 2..37) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
